$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.820.82"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.785.81"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'352.75"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").Value = "'109.01"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.552"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.602"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").Value = "'0.0838"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").Value = "'20.15"
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("D14").Value = "'7.67"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "3.221.27"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "2.811.30"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "'0.926"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "51.760.24"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'7.78"
$ws.Range("E19").Value = "  +5.13%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'13.17"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "0.0₃0966"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").Value = "'69.92"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "'266.95"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'26.16"
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +12.17%  "
$ws.Range("D29").Value = "'10.23"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'37.00"
$ws.Range("E30").Value = "  +7.97%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.22"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'6.21"
$ws.Range("E32").Value = "  +8.48%  "
$ws.Range("D33").Value = "'51.87"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("E35").Value = "  +5.22%  "
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("D39").Value = "'3.14"
$ws.Range("E39").Value = "  -2.72%  "
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").Value = "'120.65"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").Value = "'22.01"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D46").Value = "2.123.49"
$ws.Range("E46").Value = "  +2.48%  "
$ws.Range("D47").Value = "'3.31"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("E48").Value = "  +6.82%  "
$ws.Range("D49").Value = "'0.910"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "'5.44"
$ws.Range("E50").Value = "  -4.90%  "
$ws.Range("E51").Value = "  +9.05%  "

# Reset style on text-forced numeric-looking cells to remove quotePrefix marker
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
